$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.264.50"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "1.796.26"
$ws.Range("E3").Value = "  +3.14%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'338.25"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "'0.4560"
$ws.Range("E7").Value = "  +17.69%  "

$ws.Range("E8").Value = "  +11.74%  "

$ws.Range("D9").Value = "'45.13"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "'0.07593"
$ws.Range("E10").Value = "  +6.77%  "

$ws.Range("D11").Value = "'1.142"
$ws.Range("E11").Value = "  +4.41%  "

$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "'22.28"
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").Value = "'6.277"
$ws.Range("E14").Value = "  +3.80%  "

$ws.Range("D15").Value = "'7.471"
$ws.Range("E15").Value = "  +7.95%  "

$ws.Range("D16").Value = "1.798.67"
$ws.Range("E16").Value = "  +3.24%  "

$ws.Range("D17").Value = "'0.00001088"
$ws.Range("E17").Value = "  +4.53%  "

$ws.Range("D18").Value = "'0.06736"

$ws.Range("D19").Value = "'81.13"
$ws.Range("E19").Value = "  +3.43%  "

$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").Value = "'17.39"
$ws.Range("E21").Value = "  +4.86%  "

$ws.Range("D22").Value = "'6.366"
$ws.Range("E22").Value = "  +3.60%  "

$ws.Range("D23").Value = "28.254.80"
$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").Value = "'11.83"
$ws.Range("E24").Value = "  +3.59%  "

$ws.Range("D25").Value = "'2.421"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").Value = "'20.63"
$ws.Range("E26").Value = "  +5.25%  "

$ws.Range("D27").Value = "'152.01"
$ws.Range("E27").Value = "  -0.97%  "

$ws.Range("D28").Value = "'2.346"
$ws.Range("E28").Value = "  +3.71%  "

$ws.Range("D29").Value = "2.003.07"
$ws.Range("E29").Value = "  +3.29%  "

$ws.Range("D30").Value = "'132.71"
$ws.Range("E30").Value = "  +4.66%  "

$ws.Range("D31").Value = "'1.227"
$ws.Range("E31").Value = "  -2.69%  "

$ws.Range("D32").Value = "'4.029"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").Value = "'0.09460"
$ws.Range("E33").Value = "  +8.94%  "

$ws.Range("D34").Value = "'5.792"
$ws.Range("E34").Value = "  +1.36%  "

$ws.Range("D35").Value = "'0.2373"
$ws.Range("E35").Value = "  +14.59%  "

$ws.Range("D36").Value = "'0.06319"
$ws.Range("E36").Value = "  +4.97%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'12.01"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02337"
$ws.Range("E38").Value = "  +4.60%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.231"
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").Value = "'0.6581"
$ws.Range("E40").Value = "  +3.50%  "

$ws.Range("D41").Value = "'8.375"
$ws.Range("E41").Value = "  +7.68%  "

$ws.Range("D42").Value = "'1.483"
$ws.Range("E42").Value = "  -1.22%  "

$ws.Range("D43").Value = "'1.212"
$ws.Range("E43").Value = "  +2.12%  "

$ws.Range("D44").Value = "'14.20"
$ws.Range("E44").Value = "  +5.22%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").Value = "'3.849"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").Value = "'0.6083"
$ws.Range("E47").Value = "  +3.84%  "

$ws.Range("D48").Value = "'130.28"
$ws.Range("E48").Value = "  +4.19%  "

$ws.Range("D49").Value = "'2.023"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("D50").Value = "'0.07137"
$ws.Range("E50").Value = "  +3.17%  "

$ws.Range("D51").Value = "'1.162"
$ws.Range("E51").Value = "  +2.30%  "
